$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G3").Value = 'Dr. Gehan Adel, Dr. Alshimaa Atef, Administrator, Dr. Manar Montaser'
$ws.Range("G4").Value = 'Dr. Heba Mahmoud Ali, Dr. Hanan Ragab, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy'
$ws.Range("G5").Value = 'Dr. Abeer Ragab, Dr. Menna tu''Alllah Mohammad'
$ws.Range("G9").Value = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya'
$ws.Range("G10").Value = 'Dr. Basma Hamed, Dr. Amira Ibrahim'
$ws.Range("G12").Value = 'Dr. Heba Al-Sayed Mohammad, Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan'
$ws.Range("G14").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G15").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G18").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G19").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G20").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Monica, Dr. Remon, Dr. Yasmin, Dr. Nardine'
$ws.Range("G21").Value = 'Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G22").Value = 'Dr. Gehan Adel, Dr. Alshimaa Atef, Administrator, Dr. Manar Montaser'
$ws.Range("G23").Value = 'Dr. Heba Mahmoud Ali, Dr. Hanan Ragab, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy'
$ws.Range("G24").Value = 'Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G25").Value = 'Dr. Nourhan Mohammad, Dr. Yasmin Tarek'
$ws.Range("G28").Value = 'Dr. Nourhan Osama, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Yasmeena Fattoh, Dr. Dina Adel'
$ws.Range("G29").Value = 'Dr. Yasmeena Fattoh, Dr. Amira Ibrahim, Dr. Esraa Mostafa'
$ws.Range("G31").Value = 'Dr. Heba Al-Sayed Mohammad, Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan'
$ws.Range("G34").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G37").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G38").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Monica, Dr. Remon, Dr. Yasmin, Dr. Nardine'
$ws.Range("G39").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Monica, Dr. Remon, Dr. Yasmin, Dr. Nardine'
$ws.Range("G40").Value = 'Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G41").Value = 'Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Range("G42").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G43").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G47").Value = 'Dr. Nourhan Osama, Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Amira Ibrahim'
$ws.Range("G48").Value = 'Dr. Amany Raafat, Dr. Sarah Abdelmohsen, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Yasmeena Fattoh, Dr. Merna Said'
$ws.Range("G50").Value = 'Dr. Heba Al-Sayed Mohammad, Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan'
$ws.Range("G52").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G56").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G57").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Monica, Dr. Remon, Dr. Yasmin, Dr. Nardine'
$ws.Range("G58").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Monica, Dr. Remon, Dr. Yasmin, Dr. Nardine'
$ws.Range("G59").Value = 'Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G60").Value = 'Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Range("G61").Value = 'Dr. Asmaa Reda, Dr. Nahla Nagiub, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Majorelle Magdy'
$ws.Range("G63").Value = 'Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed'
$ws.Range("G66").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim'
$ws.Range("G67").Value = 'Dr. Yasmeena Fattoh, Dr. Amira Ibrahim, Dr. Esraa Mostafa'
$ws.Range("G75").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G76").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G77").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Monica, Dr. Remon, Dr. Yasmin, Dr. Nardine'
$ws.Range("G78").Value = 'Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G79").Value = 'Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Range("G80").Value = 'Dr. Asmaa Reda, Dr. Nahla Nagiub, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Majorelle Magdy'
$ws.Range("G81").Value = 'Dr. Lamiaa Ossama, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G82").Value = 'Dr. Nourhan Mohammad, Dr. Yasmin Tarek'
$ws.Range("G83").Value = 'Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed'
$ws.Range("G85").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim'
$ws.Range("G86").Value = 'Dr. Amany Raafat, Dr. Sarah Abdelmohsen, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad, Dr. Yasmeena Fattoh, Dr. Merna Said'
$ws.Range("G88").Value = 'Dr. Heba Al-Sayed Mohammad, Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan'
$ws.Range("G94").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G95").Value = 'Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G96").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Monica, Dr. Remon, Dr. Yasmin, Dr. Nardine'
$ws.Range("G97").Value = 'Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G98").Value = 'Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Range("G99").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G100").Value = 'Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Range("G101").Value = 'Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed'
$ws.Range("G104").Value = 'Dr. Nourhan Osama, Dr. Esraa Mostafa, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Amira Ibrahim'
$ws.Range("G105").Value = 'Dr. Basma Hamed, Dr. Amira Ibrahim'
$ws.Range("G113").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Range("G115").Value = 'Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Marina Sorial, Dr. Monica, Dr. Remon, Dr. Yasmin, Dr. Nardine'
